$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LeetCode Files")

foreach ($r in 30..32) {
    # Add A{r} = "x"
    $ws.Cells.Item($r, 1).Value = "x"

    # Update Difficulty (F) from 5 to 2
    $ws.Cells.Item($r, 6).Value = 2

    # Reorder concepts: G=Array, H=Binary Search, I=Matrix
    $ws.Cells.Item($r, 7).Value = "Array"
    $ws.Cells.Item($r, 8).Value = "Binary Search"
    $ws.Cells.Item($r, 9).Value = "Matrix"

    # Preserve original row height (setting values can trigger autofit on hidden rows)
    $ws.Rows.Item($r).RowHeight = 15.95
}
